$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 10218.083
$ws.Range("I11").Value = 10218.083
$ws.Range("K11").Value = 10218.083
$ws.Range("M11").Value = -10078.083
$ws.Range("H17").Value = 24476
$ws.Range("J17").Value = 25034.953
$ws.Range("L17").Value = 75104.859
$ws.Range("N17").Value = -75440.859
$ws.Range("H43").Value = 11124173
$ws.Range("I43").Value = 20002412
$ws.Range("J43").Value = 26374.75
$ws.Range("K43").Value = 20002412
$ws.Range("L43").Value = 26374.75
$ws.Range("M43").Value = -20002343
$ws.Range("N43").Value = -26512.75
$ws.Range("H69").Value = 15568.75
$ws.Range("I69").Value = 12310
$ws.Range("J69").Value = 21000
$ws.Range("K69").Value = 36930
$ws.Range("L69").Value = 63000
$ws.Range("M69").Value = -36056
$ws.Range("N69").Value = -64748
$ws.Range("H72").Value = 15568.75
$ws.Range("I72").Value = 12310
$ws.Range("J72").Value = 21000
$ws.Range("K72").Value = 110790
$ws.Range("L72").Value = 189000
$ws.Range("M72").Value = -106422
$ws.Range("N72").Value = -197736
$ws.Range("H74").Value = 9667.666999999999
$ws.Range("I74").Value = 9667.666999999999
$ws.Range("K74").Value = 9667.666999999999
$ws.Range("M74").Value = -8731.666999999999
$ws.Range("H77").Value = 9667.666999999999
$ws.Range("I77").Value = 9667.666999999999
$ws.Range("K77").Value = 48338.335
$ws.Range("M77").Value = -43658.335
$ws.Range("H100").Value = 2926.5
$ws.Range("I100").Value = 3184.9
$ws.Range("J100").Value = 1634.5
$ws.Range("K100").Value = 3184.9
$ws.Range("L100").Value = 1634.5
$ws.Range("M100").Value = -2643.9
$ws.Range("N100").Value = -2716.5
$ws.Range("H115").Value = 653.5714
$ws.Range("I115").Value = 585
$ws.Range("J115").Value = 825
$ws.Range("K115").Value = 1755
$ws.Range("L115").Value = 2475
$ws.Range("M115").Value = -188
$ws.Range("N115").Value = -5609
$ws.Range("H116").Value = 10585.467
$ws.Range("I116").Value = 14017.1
$ws.Range("K116").Value = 14017.1
$ws.Range("M116").Value = -10575.1
$ws.Range("H138").Value = 6089.9023
$ws.Range("J138").Value = 4101.364
$ws.Range("L138").Value = 12304.092
$ws.Range("N138").Value = -22584.092

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 102718
$ws.Range("I45").Value = 336164
$ws.Range("J45").Value = 2669.7144
$ws.Range("K45").Value = 336164
$ws.Range("L45").Value = 2669.7144
$ws.Range("M45").Value = -335787
$ws.Range("N45").Value = -3423.7144
$ws.Range("H61").Value = 1271371.8
$ws.Range("I61").Value = 29690.64
$ws.Range("J61").Value = 6651990
$ws.Range("K61").Value = 29690.64
$ws.Range("L61").Value = 6651990
$ws.Range("M61").Value = -29478.64
$ws.Range("N61").Value = -6652414
$ws.Range("H74").Value = 1226225.2
$ws.Range("I74").Value = 1457
$ws.Range("K74").Value = 1457
$ws.Range("M74").Value = -583
$ws.Range("H77").Value = 1226225.2
$ws.Range("I77").Value = 1457
$ws.Range("K77").Value = 7285
$ws.Range("M77").Value = -2917
$ws.Range("H97").Value = 10394.3
$ws.Range("I97").Value = 11472.111
$ws.Range("K97").Value = 11472.111
$ws.Range("M97").Value = -10976.111
$ws.Range("H105").Value = 104925.664
$ws.Range("J105").Value = 104925.664
$ws.Range("L105").Value = 104925.664
$ws.Range("N105").Value = -111913.664
$ws.Range("H132").Value = 1986
$ws.Range("I132").Value = 1882.2273
$ws.Range("K132").Value = 5646.6819
$ws.Range("M132").Value = -3116.6819
$ws.Range("H136").Value = 1271371.8
$ws.Range("I136").Value = 29690.64
$ws.Range("J136").Value = 6651990
$ws.Range("K136").Value = 89071.92
$ws.Range("L136").Value = 19955970
$ws.Range("M136").Value = -86521.92
$ws.Range("N136").Value = -19961070

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26784.615
$ws.Range("I105").Value = 34488.168
$ws.Range("J105").Value = 20181.572
$ws.Range("K105").Value = 34488.168
$ws.Range("L105").Value = 20181.572
$ws.Range("M105").Value = -32741.168
$ws.Range("N105").Value = -23675.572
$ws.Range("H107").Value = 14901.235
$ws.Range("I107").Value = 23452.2
$ws.Range("J107").Value = 2685.5715
$ws.Range("K107").Value = 23452.2
$ws.Range("L107").Value = 2685.5715
$ws.Range("M107").Value = -21532.2
$ws.Range("N107").Value = -6525.5715
$ws.Range("H134").Value = 36002868
$ws.Range("I134").Value = 2457.5715
$ws.Range("K134").Value = 7372.7145
$ws.Range("M134").Value = -4837.7145

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 805.6923
$ws.Range("I22").Value = 789.9
$ws.Range("K22").Value = 789.9
$ws.Range("M22").Value = -439.9
$ws.Range("H35").Value = 1950
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 1600
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 1600
$ws.Range("M35").Value = -1706
$ws.Range("N35").Value = -2188
$ws.Range("H58").Value = 2422.5
$ws.Range("J58").Value = 2044.5834
$ws.Range("L58").Value = 2044.5834
$ws.Range("N58").Value = -2450.5834
$ws.Range("H112").Value = 69990
$ws.Range("J112").Value = 69990
$ws.Range("L112").Value = 69990
$ws.Range("N112").Value = -72944
$ws.Range("H132").Value = 26015.928
$ws.Range("I132").Value = 38910.406
$ws.Range("K132").Value = 116731.218
$ws.Range("M132").Value = -114201.218
$ws.Range("H134").Value = 2721.4783
$ws.Range("I134").Value = 2614.4285
$ws.Range("K134").Value = 7843.2855
$ws.Range("M134").Value = -5308.2855
$ws.Range("H136").Value = 2422.5
$ws.Range("J136").Value = 2044.5834
$ws.Range("L136").Value = 6133.7502
$ws.Range("N136").Value = -11233.7502

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2194.0588
$ws.Range("I132").Value = 2066.5833
$ws.Range("K132").Value = 18599.2497
$ws.Range("M132").Value = -16069.2497
$ws.Range("H134").Value = 1435.6
$ws.Range("I134").Value = 1435.6
$ws.Range("K134").Value = 4306.799999999999
$ws.Range("M134").Value = 763.2000000000007

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H80").Value = 25811.455
$ws.Range("J80").Value = 43566.6
$ws.Range("L80").Value = 43566.6
$ws.Range("N80").Value = -45562.6
$ws.Range("H83").Value = 25811.455
$ws.Range("J83").Value = 43566.6
$ws.Range("L83").Value = 217833
$ws.Range("N83").Value = -227817
$ws.Range("H102").Value = 62501916
$ws.Range("I102").Value = 71430450
$ws.Range("K102").Value = 71430450
$ws.Range("M102").Value = -71428828
$ws.Range("H122").Value = 2573.8333
$ws.Range("J122").Value = 3100
$ws.Range("L122").Value = 9300
$ws.Range("N122").Value = -14200
$ws.Range("N44").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9843.333000000001
$ws.Range("I46").Value = 14335.25
$ws.Range("K46").Value = 14335.25
$ws.Range("M46").Value = -14147.25
$ws.Range("H136").Value = 71255.13
$ws.Range("I136").Value = 146807.28
$ws.Range("K136").Value = 440421.84
$ws.Range("M136").Value = -437871.84

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7319.6665
$ws.Range("J41").Value = 7431.857
$ws.Range("L41").Value = 7431.857
$ws.Range("N41").Value = -8211.857
$ws.Range("H100").Value = 1027.8334
$ws.Range("I100").Value = 993.4
$ws.Range("K100").Value = 1986.8
$ws.Range("M100").Value = -1445.8
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H114").Value = 27500
$ws.Range("J114").Value = 27500
$ws.Range("L114").Value = 27500
$ws.Range("N114").Value = -36178
$ws.Range("H122").Value = 1661.9286
$ws.Range("I122").Value = 1800.8182
$ws.Range("K122").Value = 5402.4546
$ws.Range("M122").Value = -2952.4546
$ws.Range("H126").Value = 2305.75
$ws.Range("I126").Value = 2216.2727
$ws.Range("K126").Value = 6648.8181
$ws.Range("M126").Value = -4178.8181
$ws.Range("H136").Value = 39602.81
$ws.Range("I136").Value = 53322.895
$ws.Range("J136").Value = 2362.5715
$ws.Range("K136").Value = 159968.685
$ws.Range("L136").Value = 7087.7145
$ws.Range("M136").Value = -157418.685
$ws.Range("N136").Value = -12187.7145
